$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1916
$ws1.Range("F5").Value = 306
$ws1.Range("F6").Value = 44
$ws1.Range("F7").Value = 90
$ws1.Range("F8").Value = 2040
$ws1.Range("F9").Value = 10340
$ws1.Range("F12").Value = 266
$ws1.Range("F13").Value = 197
$ws1.Range("F15").Value = 7193
$ws1.Range("F17").Value = 679
$ws1.Range("F18").Value = 118
$ws1.Range("F20").Value = 259

# Sheet "全部类型" (sheet4) - update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1916
$ws4.Range("F5").Value = 306
$ws4.Range("F6").Value = 44
$ws4.Range("F8").Value = 90
$ws4.Range("F9").Value = 2040
$ws4.Range("F12").Value = 10340
$ws4.Range("F15").Value = 266
$ws4.Range("F16").Value = 197
$ws4.Range("F18").Value = 7193
$ws4.Range("F20").Value = 679
$ws4.Range("F21").Value = 118
$ws4.Range("F23").Value = 259
